# Weekly update: insert two new price rows at the top of the data block
# (right before the existing row 828), pushing all subsequent rows down
# by two. The two new rows carry a new reporting date (serial 45142 =
# 2023-08-04) for the "Pintón" and "Primera Pintón" quality grades.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 828; everything currently
# at/after row 828 shifts down to row 830+.
$ws.Rows("828:829").Insert()

# ---- New row 828 : Pintón ----
$ws.Cells.Item(828, 1).Value = 4
$ws.Cells.Item(828, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(828, 3).Value = "Los Lagos"
$ws.Cells.Item(828, 4).Value = 45142
$ws.Cells.Item(828, 5).Value = 10
$ws.Cells.Item(828, 6).Value = "Fruta"
$ws.Cells.Item(828, 7).Value = 100108
$ws.Cells.Item(828, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(828, 9).Value = 100108006
$ws.Cells.Item(828, 10).Value = "Plátano"
$ws.Cells.Item(828, 11).Value = "Sin especificar"
$ws.Cells.Item(828, 12).Value = "Pintón"
$ws.Cells.Item(828, 13).Value = 600
$ws.Cells.Item(828, 14).Value = 17000
$ws.Cells.Item(828, 15).Value = 17000
$ws.Cells.Item(828, 16).Value = 17000
$ws.Cells.Item(828, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(828, 18).Value = "Ecuador"
$ws.Cells.Item(828, 19).Value = 850
$ws.Cells.Item(828, 20).Value = 20

# ---- New row 829 : Primera Pintón ----
$ws.Cells.Item(829, 1).Value = 4
$ws.Cells.Item(829, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(829, 3).Value = "Los Lagos"
$ws.Cells.Item(829, 4).Value = 45142
$ws.Cells.Item(829, 5).Value = 10
$ws.Cells.Item(829, 6).Value = "Fruta"
$ws.Cells.Item(829, 7).Value = 100108
$ws.Cells.Item(829, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(829, 9).Value = 100108006
$ws.Cells.Item(829, 10).Value = "Plátano"
$ws.Cells.Item(829, 11).Value = "Sin especificar"
$ws.Cells.Item(829, 12).Value = "Primera Pintón"
$ws.Cells.Item(829, 13).Value = 600
$ws.Cells.Item(829, 14).Value = 19000
$ws.Cells.Item(829, 15).Value = 19000
$ws.Cells.Item(829, 16).Value = 19000
$ws.Cells.Item(829, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(829, 18).Value = "Ecuador"
$ws.Cells.Item(829, 19).Value = 950
$ws.Cells.Item(829, 20).Value = 20

# Make sure the D column on the two new rows keeps the date style (s="2")
# used by every other row in that column.
$ws.Range("D828:D829").NumberFormat = $ws.Range("D830").NumberFormat
